# Updated cryptos list - refresh price/volume(1h) figures and re-rank a
# few coins (Hedera/Monero/EthereumClassic and Kaspa/LidoDAOToken swapped).
#
# Note: several "Price" values look numeric (e.g. "321.37") but must stay
# stored as text, exactly like the rest of the column. Assigning such a
# string straight to .Value would make Excel coerce it into a real number,
# so for those cells we briefly force a text NumberFormat, set the value,
# then reset the cell style back to "Normal" so no visible formatting
# change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = $value
}

function Set-NumericLookingTextValue($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" '44.527.93'
Set-TextValue "E2" '  +3.81%  '

# Row 3 - Ethereum
Set-TextValue "D3" '2.284.05'
Set-TextValue "E3" '  +2.77%  '

# Row 4 - TetherUSD
Set-TextValue "E4" '  +0.02%  '

# Row 5 - BNB
Set-NumericLookingTextValue "D5" '321.37'

# Row 6 - Solana
Set-NumericLookingTextValue "D6" '106.34'
Set-TextValue "E6" '  +6.15%  '

# Row 7 - XRP
Set-TextValue "E7" '  +0.13%  '

# Row 8 - USDC
Set-TextValue "E8" '  +0.03%  '

# Row 9 - Cardano
Set-TextValue "E9" '  +1.68%  '

# Row 10 - Avalanche
Set-NumericLookingTextValue "D10" '38.87'
Set-TextValue "E10" '  +4.14%  '

# Row 11 - Dogecoin
Set-NumericLookingTextValue "D11" '0.0845'
Set-TextValue "E11" '  +1.71%  '

# Row 12 - Polkadot
Set-NumericLookingTextValue "D12" '7.93'
Set-TextValue "E12" '  +1.45%  '

# Row 13 - TRON
Set-NumericLookingTextValue "D13" '0.107'
Set-TextValue "E13" '  +1.10%  '

# Row 14 - Polygon
Set-TextValue "E14" '  +2.45%  '

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" '2.630.17'
Set-TextValue "E15" '  +2.71%  '

# Row 16 - Chainlink
Set-TextValue "E16" '  +2.58%  '

# Row 17 - WrappedEther
Set-TextValue "D17" '2.340.28'
Set-TextValue "E17" '  +5.44%  '

# Row 18 - WrappedBTC
Set-TextValue "D18" '44.366.46'
Set-TextValue "E18" '  +3.66%  '

# Row 19 - InternetComputer(DFINITY)
Set-TextValue "E19" '  -5.90%  '

# Row 20 - ShibaInu
Set-TextValue "E20" '  +4.24%  '

# Row 21 - Uniswap
Set-TextValue "E21" '  +1.31%  '

# Row 22 - Litecoin
Set-NumericLookingTextValue "D22" '66.50'
Set-TextValue "E22" '  +1.45%  '

# Row 23 - PancakeSwap
Set-TextValue "E23" '  +1.81%  '

# Row 24 - BitcoinCash
Set-NumericLookingTextValue "D24" '239.90'
Set-TextValue "E24" '  +1.37%  '

# Row 25 - ImmutableX
Set-TextValue "E25" '  +3.22%  '

# Row 26 - Dai
Set-TextValue "E26" '  -0.53%  '

# Row 27 - Cosmos
Set-NumericLookingTextValue "D27" '10.23'
Set-TextValue "E27" '  +1.12%  '

# Row 28 - InjectiveProtocol
Set-NumericLookingTextValue "D28" '38.56'
Set-TextValue "E28" '  +12.45%  '

# Row 29 - Toncoin
Set-TextValue "E29" '  +0.44%  '

# Row 30 - Filecoin
Set-TextValue "E30" '  +2.24%  '

# Row 31 - was Hedera, now Monero
Set-TextValue "B31" 'Monero'
Set-TextValue "C31" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-NumericLookingTextValue "D31" '163.11'
Set-TextValue "E31" '  +4.03%  '

# Row 32 - was Monero, now EthereumClassic
Set-TextValue "B32" 'EthereumClassic'
Set-TextValue "C32" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-NumericLookingTextValue "D32" '20.68'
Set-TextValue "E32" '  +0.76%  '

# Row 33 - was EthereumClassic, now Hedera
Set-TextValue "B33" 'Hedera'
Set-TextValue "C33" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-NumericLookingTextValue "D33" '0.0890'
Set-TextValue "E33" '  -2.23%  '

# Row 34 - WEMIXToken
Set-TextValue "E34" '  -1.35%  '

# Row 35 - ARBITRUM
Set-TextValue "E35" '  +5.08%  '

# Row 36 - was LidoDAOToken, now Kaspa
Set-TextValue "B36" 'Kaspa'
Set-TextValue "C36" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-NumericLookingTextValue "D36" '0.118'
Set-TextValue "E36" '  +11.93%  '

# Row 37 - was Kaspa, now LidoDAOToken
Set-TextValue "B37" 'LidoDAOToken'
Set-TextValue "C37" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-NumericLookingTextValue "D37" '3.22'
Set-TextValue "E37" '  +0.53%  '

# Row 38 - Stellar
Set-TextValue "E38" '  -0.45%  '

# Row 39 - NEARProtocol
Set-NumericLookingTextValue "D39" '3.98'
Set-TextValue "E39" '  -0.02%  '

# Row 40 - RenderToken
Set-TextValue "E40" '  +0.18%  '

# Row 41 - VeChain
Set-TextValue "E41" '  +1.22%  '

# Row 42 - Celestia
Set-NumericLookingTextValue "D42" '15.51'
Set-TextValue "E42" '  +24.55%  '

# Row 43 - FirstDigitalUSD
Set-TextValue "E43" '  +0.12%  '

# Row 44 - Maker
Set-TextValue "D44" '1.777.41'
Set-TextValue "E44" '  -8.76%  '

# Row 45 - Algorand
Set-TextValue "E45" '  +0.23%  '

# Row 46 - BitcoinSV
Set-NumericLookingTextValue "D46" '87.02'
Set-TextValue "E46" '  -2.07%  '

# Row 47 - THORChain
Set-TextValue "E47" '  +2.58%  '

# Row 48 - MultiversX
Set-NumericLookingTextValue "D48" '60.59'
Set-TextValue "E48" '  -0.06%  '

# Row 49 - Stacks
Set-TextValue "E49" '  +5.87%  '

# Row 50 - ordi
Set-NumericLookingTextValue "D50" '74.74'
Set-TextValue "E50" '  -2.20%  '

# Row 51 - FraxShare
Set-TextValue "E51" '  +1.27%  '
